# "finish validation analysis and add noise robustness"
#
# Adds a new "validation" worksheet (after "immune_cells") that mirrors the
# "immune_cells" value/color table but swaps the last two cell-type rows for
# "Neutrophil" and "T cell" (a noise-robustness check against the existing
# cell-type palette). The new sheet becomes the active tab.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the tab strip -------------------------
$immuneCells = $wb.Worksheets.Item("immune_cells")
$ws = $wb.Worksheets.Add($null, $immuneCells)
$ws.Name = "validation"

# --- same "value"/"color" table as immune_cells, with row 8/9 replaced -----
$ws.Range("A1").Value = "value"
$ws.Range("B1").Value = "color"

$ws.Range("A2").Value = "B cell"
$ws.Range("B2").Value = "#1b9e77"

$ws.Range("A3").Value = "Dendritic cell"
$ws.Range("B3").Value = "#d95f02"

$ws.Range("A4").Value = "Macrophage/Monocyte"
$ws.Range("B4").Value = "#7570b3"

$ws.Range("A5").Value = "NK cell"
$ws.Range("B5").Value = "#e7298a"

$ws.Range("A6").Value = "T cell CD4+"
$ws.Range("B6").Value = "#66a61e"

$ws.Range("A7").Value = "T cell CD8+"
$ws.Range("B7").Value = "#e6ab02"

$ws.Range("A8").Value = "Neutrophil"
$ws.Range("B8").Value = "#a6761d"

$ws.Range("A9").Value = "T cell"
$ws.Range("B9").Value = "#666666"

# column A is a long text column on every one of these sheets - match it
$ws.Columns.Item(1).ColumnWidth = 25.33

# leave the new sheet active, selection parked one row below the data
$ws.Activate()
$ws.Range("A10").Select() | Out-Null
